$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# New data rows (NC, Paterno, Materno, Nombres, Nombre_Largo, Grupo, Reprobadas)
$data = @(
    @(21330051920005, "CASTILLO",  "GARCIA",    "KEVIN ISAAC",       "LÓGICA", "1AV", 6),
    @(21330051920006, "CARRERA",   "HERNANDEZ", "EDGAR FLORENCIO",   "LÓGICA", "1AV", 6),
    @(21330051920017, "MARTINEZ",  "XOTLANIHUA","YAIR",               "LÓGICA", "1AV", 6),
    @(21330051920379, "MORALES",   "SANCHEZ",   "MIGUEL",             "LÓGICA", "1AV", 6),
    @(21330051920067, "VEGA",      "ROJAS",     "HANIA ZARETH",       "LÓGICA", "1BV", 6),
    @(21330051920078, "GARCIA",    "GONZALEZ",  "MIROSLAVA",          "LÓGICA", "1CV", 6),
    @(21330051920101, "ROMERO",    "RAMIREZ",   "CITLALI ESPERANZA",  "LÓGICA", "1CV", 6),
    @(21330051920355, "GARCIA",    "GASPAR",    "MARIA VALERIA",      "LÓGICA", "1AV", 6),
    @(21330051920126, "MARIN",     "RODRIGUEZ", "ALEJANDRO",          "LÓGICA", "1DV", 6),
    @(20330051920275, "MONTIEL",   "FLORES",    "IKER XAVIER",        "LÓGICA", "1DV", 2)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $row = $row + 1
}
